$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Plan Start" (column E) values for rows 24 & 25 ---
$ws.Cells.Item(24, 5).Value = 50
$ws.Cells.Item(25, 5).Value = 50

# --- Row 26 stays "3.6 Implement Food Wars Feature" (unchanged text, kept for clarity) ---
$ws.Cells.Item(26, 2).Value = "3.6 Implement Food Wars Feature"

# --- Row 27: was "3.7 Conduct Testing" -> now "3.7 Executive Summary" ---
$ws.Cells.Item(27, 2).Value = "3.7 Executive Summary"
$ws.Cells.Item(27, 3).Value = 48
$ws.Cells.Item(27, 4).Value = 3
$ws.Cells.Item(27, 7).ClearContents()

# --- Row 28: new activity "3.8 Unit Testing and Report" (was "4.1 Track Progress").
#     D28 switches from a text ("Ongoing") cell to a plain-number cell, so repaint its
#     format from a numeric neighbour (C28) before writing the new number into it. ---
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Cells.Item(28, 2).Value = "3.8 Unit Testing and Report"
$ws.Cells.Item(28, 3).Value = 51
$ws.Cells.Item(28, 4).Value = 2
$ws.Cells.Item(28, 5).ClearContents()
$ws.Cells.Item(28, 7).ClearContents()

# --- Row 29: new activity "3.9 Coverage Testing and Report" (was "4.2 Control Scope") ---
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Cells.Item(29, 2).Value = "3.9 Coverage Testing and Report"
$ws.Cells.Item(29, 3).Value = 53
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).ClearContents()
$ws.Cells.Item(29, 7).ClearContents()

# --- Row 30: now "4.1 Track Progress" (was "4.3 Control Quality") ---
$ws.Cells.Item(30, 2).Value = "4.1 Track Progress"
$ws.Cells.Item(30, 3).Value = 1
$ws.Cells.Item(30, 4).Value = "Ongoing"
$ws.Cells.Item(30, 5).Value = 46

# --- Row 31: now "4.2 Control Scope" (was "4.4 Monitor Bugs") ---
$ws.Cells.Item(31, 2).Value = "4.2 Control Scope"
$ws.Cells.Item(31, 3).Value = 30
$ws.Cells.Item(31, 4).Value = "Ongoing"
$ws.Cells.Item(31, 5).Value = 46

# --- Row 32: now "4.3 Control Quality" (was "5.1 Review and Complete Work").
#     D32 switches from a plain number to a text ("Ongoing") cell, so repaint its format
#     from an existing "Ongoing" cell (D30) before writing the text into it. ---
$ws.Range("D30").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Cells.Item(32, 2).Value = "4.3 Control Quality"
$ws.Cells.Item(32, 3).Value = 30
$ws.Cells.Item(32, 4).Value = "Ongoing"
$ws.Cells.Item(32, 5).Value = 46

# --- Row 33: now "4.4 Monitor Bugs" (was "5.2 Performance Review") ---
$ws.Range("D30").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Cells.Item(33, 2).Value = "4.4 Monitor Bugs"
$ws.Cells.Item(33, 3).Value = 30
$ws.Cells.Item(33, 4).Value = "Ongoing"
$ws.Cells.Item(33, 5).Value = 46

# --- Row 34: now "5.1 Review and Complete Work" (was "5.3 Verify Completion of Work") ---
$ws.Cells.Item(34, 2).Value = "5.1 Review and Complete Work"
$ws.Cells.Item(34, 3).Value = 55
$ws.Cells.Item(34, 4).Value = 1

# --- Row 35: new row "5.2 Performance Review" (was a blank spacer row). Bring over the
#     same cell formatting (styles 23/5/5/5/5/22) used by the row above, then set values. ---
$ws.Range("B34:G34").Copy()
$ws.Range("B35:G35").PasteSpecial(-4122)
$ws.Cells.Item(35, 2).Value = "5.2 Performance Review"
$ws.Cells.Item(35, 3).Value = 56
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 7).Value = 0

# --- Row 36: new row "5.3 Verify Completion of Work" (was a blank spacer row) ---
$ws.Range("B34:G34").Copy()
$ws.Range("B36:G36").PasteSpecial(-4122)
$ws.Cells.Item(36, 2).Value = "5.3 Verify Completion of Work"
$ws.Cells.Item(36, 3).Value = 57
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 7).Value = 0

$excel.CutCopyMode = 0

# --- Two new "decorative" blank rows (41 & 42), matching the existing blank rows 37-40 ---
$ws.Rows.Item(41).RowHeight = 30
$ws.Rows.Item(42).RowHeight = 30
$ws.Range("B41:AA42").Value = 0
$ws.Range("B41:AA42").ClearContents()

# --- Conditional formatting now covers through row 36 instead of row 34 ---
$fcs = $ws.Range("H5:BS34").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("H5:BS36"))
}

# --- Selection / view state ---
$ws.Range("A16").Select()
$ws.Range("S27").Select()
